$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

# Column D holds a date value; copy the style from the row above so the
# numeric date format (s="2") is preserved, then set the date value.
$ws.Range("D4").Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0).Date

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112044
$ws.Cells.Item($row, 7).Value = "Perejil"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 600
$ws.Cells.Item($row, 12).Value = 650
$ws.Cells.Item($row, 13).Value = 625
$ws.Cells.Item($row, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item($row, 15).Value = "Región del Maule"
$ws.Cells.Item($row, 16).Value = 625
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
